$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2:G2) currently holds empty-string inline-text cells; clear them
# out entirely so the row becomes a bare, cell-less row.
$ws.Range("A2:G2").ClearContents()
# Touch the row itself (toggle + restore a default row property) so the
# now-empty row 2 is still emitted as a bare <row r="2"/> element instead
# of being dropped from sheetData entirely.
$ws.Rows(2).Hidden = $true
$ws.Rows(2).Hidden = $false

# Append two new rows (109, 110) of "test create" sample data.
# Column A holds a number-looking value that must stay text (leading
# apostrophe forces Excel to store it as text, matching the rest of column A).
$ws.Range("A109").Value = "'107"
$ws.Range("B109").Value = "test create"
$ws.Range("C109").Value = "omar"
$ws.Range("D109").Value = "egypt"
$ws.Range("E109").Value = "https://ar.wikipedia.org/wiki/%D8%B5%D9%86%D8%B9_%D8%A7%D9%84%D9%84%D9%87_%D8%A5%D8%A8%D8%B1%D8%A7%D9%87%D9%8A%D9%85"
$ws.Range("F109").Value = "https://ar.wikipedia.org/wiki/%D9%85%D8%B5%D8%B1"
$ws.Range("G109").Value = "https://ar.wikipedia.org/wiki/%D9%85%D8%B5%D8%B1"

$ws.Range("A110").Value = "'108"
$ws.Range("B110").Value = "test create"
$ws.Range("C110").Value = "omar"
$ws.Range("D110").Value = "egypt"
$ws.Range("E110").Value = "https://ar.wikipedia.org/wiki/%D8%B5%D9%86%D8%B9_%D8%A7%D9%84%D9%84%D9%87_%D8%A5%D8%A8%D8%B1%D8%A7%D9%87%D9%8A%D9%85"
$ws.Range("F110").Value = "https://ar.wikipedia.org/wiki/%D9%85%D8%B5%D8%B1"
$ws.Range("G110").Value = "https://ar.wikipedia.org/wiki/%D9%85%D8%B5%D8%B1"
